$p = $ppt.ActivePresentation

# --- Duplicate the "Monday, January 27" slide (slide 6) to create the new
# "Wednesday, January 29" slide. PowerPoint places the duplicate right
# after the source slide, i.e. at index 7 (the new last slide). ---
$s6 = $p.Slides.Item(6)
$dup = $s6.Duplicate()
$s7 = $p.Slides.Item(7)

# --- Title placeholder ---
$title = $s7.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Wednesday, January 29"

# --- Body placeholder ---
$body = $s7.Shapes.Item(2).TextFrame.TextRange

$p1 = "First Quizzes are today!"
$p2 = "Gradescope appears to be working well so far"
$p3 = "Board game: Changed description of input slightly."
$p4 = "Will change the due dates on those very very soon. "
$p5 = "Please join the class Piazza asap. Tas are posting common questions there."
$p6 = "Today we will continue and possibly finish discussing advanced graphs"

$body.Text = $p1 + "`r" + $p2 + "`r" + $p3 + "`r" + $p4 + "`r" + $p5 + "`r" + $p6

# Paragraph 3 & 4 are sub-bullets (indent level 2 == lvl="1")
$body.Lines(3,1).IndentLevel = 2
$body.Lines(4,1).IndentLevel = 2

# Paragraph 1: "First " + bold/italic "Quizzes" + " are today!"
$off = 1
$r = $body.Characters($off + ("First ").Length, ("Quizzes").Length)
$r.Font.Bold = -1
$r.Font.Italic = -1

# Paragraph 2: bold/italic "Gradescope" + " appears to be working well so far"
$off = $off + $p1.Length + 1
$r = $body.Characters($off, ("Gradescope").Length)
$r.Font.Bold = -1
$r.Font.Italic = -1

# Paragraph 5: "Please join the class " + bold/italic "Piazza" + " asap. " + "Tas" + " are posting..."
$off = $off + $p2.Length + 1 + $p3.Length + 1 + $p4.Length + 1
$r = $body.Characters($off + ("Please join the class ").Length, ("Piazza").Length)
$r.Font.Bold = -1
$r.Font.Italic = -1

# Paragraph 6: "... discussing " + bold/italic "advanced graphs"
$off = $off + $p5.Length + 1
$r = $body.Characters($off + ("Today we will continue and possibly finish discussing ").Length, ("advanced graphs").Length)
$r.Font.Bold = -1
$r.Font.Italic = -1

Write-Output ("Slide count: " + $p.Slides.Count)
